$wb = $excel.ActiveWorkbook

# --- Update selection on sheet "2_" (selection sqref="A1:D5" -> activeCell="B1" sqref="B1:D5") ---
$ws2 = $wb.Worksheets.Item("2_")
$ws2.Range("B1:D5").Select()

# --- Update selection on sheet "17_" (selection sqref="A1:C3" -> activeCell="E3" sqref="E3") ---
$ws17 = $wb.Worksheets.Item("17_")
$ws17.Range("E3").Select()

# --- Add the new sheet "19_" at the very end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "19_"

# Column widths (closest achievable values given engine's column-width rounding granularity)
$newSheet.Columns.Item(1).ColumnWidth = 34.166666666666664
$newSheet.Columns.Item(3).ColumnWidth = 32.666666666666664
$newSheet.Columns.Item(4).ColumnWidth = 53.666666666666664

# Row heights
$newSheet.Rows.Item(1).RowHeight = 150
$newSheet.Rows.Item(2).RowHeight = 30
$newSheet.Rows.Item(3).RowHeight = 30
$newSheet.Rows.Item(4).RowHeight = 45
$newSheet.Rows.Item(5).RowHeight = 45
$newSheet.Rows.Item(6).RowHeight = 45

# --- Values ---
# Row 1
$newSheet.Range("A1").Value = "In class, we will attempt to simulate a system governed by the differential equation above.   Assume that we will keep the circuit details (size of the resistor, etc…) the same and vary the input signal;  also assume that we want to know much of the input signal makes it through the filter, match the variables in the equation to the type of variable. "
$newSheet.Range("B1").Value = "Correct order of definitions"
$newSheet.Range("C1").Value = "Definitions"

# Row 2
$newSheet.Range("A2").Value = "R"
$newSheet.Range("B2").Value = "C"
$newSheet.Range("C2").Value = "State Variable"
$newSheet.Range("D2").Value = "The size of the resistor R is a physical part of the circuit: here we are keeping that constant."

# Row 3
$newSheet.Range("A3").Value = "C"
$newSheet.Range("B3").Value = "C"
$newSheet.Range("C3").Value = "Metric"
$newSheet.Range("D3").Value = "The size of the capacitor C is a physical part of the circuit: here we are keeping that constant."

# Row 4
$newSheet.Range("A4").Value = "Frequency and amplitude of V_in"
$newSheet.Range("B4").Value = "D"
$newSheet.Range("C4").Value = "Parameter"
$newSheet.Range("D4").Value = "The amplitude and frequency of the input signal is an independent variable: we are interested in how the filter reacts to different input signals"

# Row 5
$newSheet.Range("A5").Value = "Value of V_out at a given moment"
$newSheet.Range("B5").Value = "A"
$newSheet.Range("C5").Value = "Independent variable"
$newSheet.Range("D5").Value = "Like T in our coffee problem, V_out is our state variable: at any given moment, the ""state of the system"" is described by this value"

# Row 6 (note: B6 uses the plain wrap style, not the centered style used by B1:B5)
$newSheet.Range("A6").Value = "Amplitude of V_out"
$newSheet.Range("B6").Value = "B"
$newSheet.Range("D6").Value = "The filter will reduce the amplitude of high frequency signals; this amplitude tells us which signals have been filtered and which made it through the filter."

# --- Styles ---
# Column A rows 1:6 - plain wrap
$newSheet.Range("A1:A6").WrapText = $true

# Column B rows 1:5 - centered + wrapped (the "metric" labels column)
$bCentered = $newSheet.Range("B1:B5")
$bCentered.WrapText = $true
$bCentered.HorizontalAlignment = -4108
$bCentered.VerticalAlignment = -4108

# B6 - plain wrap (different from B1:B5)
$newSheet.Range("B6").WrapText = $true

# Column C rows 1:6 - plain wrap (also materializes the empty C6 cell)
$newSheet.Range("C1:C6").WrapText = $true

# Column D rows 1:6 - plain wrap (also materializes the empty D1 cell)
$newSheet.Range("D1:D6").WrapText = $true

# Rows 7:14 across columns A:D - empty, plain wrap style
$newSheet.Range("A7:D14").WrapText = $true

# --- Finally select C1 on the new sheet, matching the target selection state ---
$newSheet.Range("C1").Select()
